$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# Fill in the two new strings for the previously empty cells A36/A37
$ws.Range("A36").Value = "Выберите реакцию "
$ws.Range("A37").Value = "Добавление канала было отменено"

# Row 37 grows taller to fit the wrapped text
$ws.Rows.Item(37).RowHeight = 30

# Update the saved view/selection state to match the new scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("A37").Select()
